# W9S3 materials finish-up edit.
# 1) Bump the cached "datetimeFigureOut" date fields from 20/3/2023 to
#    21/3/2023 across every slide layout, the slide master and the notes
#    master (these hold the cached field text shown in the diff).
# 2) Fix the "Litteral"/"litterals" typo on slide 24 ("Lexical Analysis")
#    in the "Content Placeholder 2" shape, matching the new run layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder refresh: 20/3/2023 -> 21/3/2023
# ---------------------------------------------------------------------

function Update-DateField($shape) {
    if ($shape.HasTextFrame -eq -1) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "20/3/2023") {
            $tr.Characters(1, $tr.Length).Text = "21/3/2023"
        }
    }
}

# Slide master's own Date Placeholder.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateField $master.Shapes.Item($i)
}

# Every slide layout's Date Placeholder.
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $lay = $layouts.Item($L)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        Update-DateField $lay.Shapes.Item($i)
    }
}

# Notes master's Date Placeholder.
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    Update-DateField $notesMaster.Shapes.Item($i)
}

# ---------------------------------------------------------------------
# 2) Slide 24 ("Lexical Analysis") text fixes
# ---------------------------------------------------------------------

$s24 = $p.Slides.Item(24)
$sh = $s24.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# Paragraph 1: "Litteral: " -> "Literal: " (merge the two leading runs).
$para1 = $tr.Paragraphs(1)
$para1.Characters(1, 10).Text = "Literal: "

# Paragraph 2: "Examples of litterals in the C programming language: "
#            -> "Examples of literals in the C programming language: "
#   with the run break moved so "of literals " becomes its own run.
$para2 = $tr.Paragraphs(2)
$para2.Characters(10, 13).Text = "of literals "
